# Result DF Display.xlsx - "added warnings and alerts and more logging to columns algorithm"
#
# 1) Column H header (H1) text changes from "fess" to "0010".
# 2) Column H (rows 2-31) formulas change from
#       =SUBSTITUTE(Dn:Dn,"@","")
#    to
#       =IF(YEAR(Gn:Gn)<2030, An:An, "")
#    i.e. flag/alert rows whose "mydate" (column G) year is before 2030 by
#    echoing back the row id (column A), otherwise leave blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the header text in H1 ---
# A leading apostrophe forces Excel to store the numeric-looking text
# "0010" as literal text (quote-prefixed) instead of auto-converting it
# to the number 10.
$ws.Range("H1").Value = "'0010"

# --- 2) Update the formulas in H2:H31 ---
$lastRow = 31
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Range("H$row").Formula = "=IF(YEAR(G$row`:G$row)<2030, A$row`:A$row, `"`")"
}
